$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-17 06:03:49"
$wsZhCn.Range("G4").Value = "2016-02-17 06:04:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-17 06:03:59"
$wsDeDe.Range("G4").Value = "2016-02-17 06:04:54"
